# Update the "想去人数" (F column) counts for various events across the
# four sheets of the 广州-漫展信息 workbook. Values below were derived by
# comparing the pre-edit and post-edit cell contents.

$wb = $excel.ActiveWorkbook

# 展览 (Sheet 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 16
$ws1.Range("F4").Value = 1346
$ws1.Range("F7").Value = 3944
$ws1.Range("F8").Value = 245
$ws1.Range("F9").Value = 797
$ws1.Range("F10").Value = 2397
$ws1.Range("F11").Value = 373
$ws1.Range("F12").Value = 55
$ws1.Range("F13").Value = 241
$ws1.Range("F15").Value = 217
$ws1.Range("F17").Value = 3789
$ws1.Range("F18").Value = 326
$ws1.Range("F20").Value = 52
$ws1.Range("F23").Value = 58
$ws1.Range("F24").Value = 286

# 演出 (Sheet 2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 34
$ws2.Range("F22").Value = 86

# 本地生活 (Sheet 3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6403
$ws3.Range("F5").Value = 353
$ws3.Range("F7").Value = 1

# 全部类型 (Sheet 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6403
$ws4.Range("F5").Value = 353
$ws4.Range("F7").Value = 16
$ws4.Range("F10").Value = 1346
$ws4.Range("F12").Value = 34
$ws4.Range("F17").Value = 3944
$ws4.Range("F19").Value = 245
$ws4.Range("F22").Value = 797
$ws4.Range("F23").Value = 2397
$ws4.Range("F24").Value = 373
$ws4.Range("F25").Value = 55
$ws4.Range("F27").Value = 241
$ws4.Range("F29").Value = 217
$ws4.Range("F33").Value = 326
$ws4.Range("F37").Value = 52
$ws4.Range("F40").Value = 58
$ws4.Range("F47").Value = 1
$ws4.Range("F48").Value = 86
$ws4.Range("F49").Value = 286
